# Generate Report for handoff
#
# The "Latest Handoff Datetime" (column D) for the e599cfc2-... file (row 5)
# gets a fresh handoff timestamp on both the zh-cn and de-de report sheets.
$wb = $excel.ActiveWorkbook

$zhSheet = $wb.Worksheets.Item("zh-cn")
$deSheet = $wb.Worksheets.Item("de-de")

$zhSheet.Range("D5").Value = "2016-01-17 15:57:55"
$deSheet.Range("D5").Value = "2016-01-17 15:58:05"
